# 9.3.2.xlsx gh-pages deploy: add the 2022 data column (R) to the table,
# mirroring the formatting already used for the 2021 column (Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 4): new year label 2022, formatted like Q4.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R4").Value = 2022

# Data row 5: new figure, formatted like Q5.
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R5").Value = 8.6821914120339212

# Data row 6: new figure, formatted like Q6.
$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R6").Value = 12.221423436376707

# Matches the author's recorded selection after the edit.
$ws.Range("S4").Select()
